$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows 270-301 (2021-05-28 .. 2021-06-28), matching the source
# report's style for column A (date format carried from the last existing
# row) and plain numeric values for B/C/D.
$ws.Range("A269").Copy()
$ws.Range("A270:A301").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @(270, 44344, 0, 4, 87.24100327153762),
    @(271, 44345, 1, 4, 87.24100327153762),
    @(272, 44346, 0, 2, 43.62050163576881),
    @(273, 44347, 1, 3, 65.43075245365321),
    @(274, 44348, 0, 3, 65.43075245365321),
    @(275, 44349, 0, 3, 65.43075245365321),
    @(276, 44350, 0, 2, 43.62050163576881),
    @(277, 44351, 0, 2, 43.62050163576881),
    @(278, 44352, 0, 1, 21.81025081788441),
    @(279, 44353, 0, 1, 21.81025081788441),
    @(280, 44354, 0, 0, 0),
    @(281, 44355, 0, 0, 0),
    @(282, 44356, 0, 0, 0),
    @(283, 44357, 0, 0, 0),
    @(284, 44358, 0, 0, 0),
    @(285, 44359, 0, 0, 0),
    @(286, 44360, 0, 0, 0),
    @(287, 44361, 0, 0, 0),
    @(288, 44362, 0, 0, 0),
    @(289, 44363, 0, 0, 0),
    @(290, 44364, 0, 0, 0),
    @(291, 44365, 0, 0, 0),
    @(292, 44366, 0, 0, 0),
    @(293, 44367, 0, 0, 0),
    @(294, 44368, 0, 0, 0),
    @(295, 44369, 0, 0, 0),
    @(296, 44370, 0, 0, 0),
    @(297, 44371, 0, 0, 0),
    @(298, 44372, 0, 0, 0),
    @(299, 44373, 0, 0, 0),
    @(300, 44374, 0, 0, 0),
    @(301, 44375, 0, 0, 0)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
}

